$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 49: new work-log entry
$ws.Range("A49").Value = (Get-Date -Year 2018 -Month 4 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B49").Value = "Travail en groupe"
$ws.Range("C49").Value = 1.5

# Row 50: new work-log entry (hours left blank)
$ws.Range("A50").Value = (Get-Date -Year 2018 -Month 4 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B50").Value = "Travail en groupe via chat vocal"

# Update the active selection to reflect where the user ended up
$ws.Range("F46").Select() | Out-Null
